$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 365, shifting existing rows 365-414 down to 366-415.
$ws.Rows(365).Insert()

# Populate the newly inserted row 365 with the new weekly price entry.
$ws.Cells.Item(365, 1).Value = 4
$ws.Cells.Item(365, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(365, 3).Value = "Los Lagos"
$ws.Cells.Item(365, 4).Value = 45127
$ws.Cells.Item(365, 5).Value = 10
$ws.Cells.Item(365, 6).Value = "Fruta"
$ws.Cells.Item(365, 7).Value = 100108
$ws.Cells.Item(365, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(365, 9).Value = 100108002
$ws.Cells.Item(365, 10).Value = "Mango"
$ws.Cells.Item(365, 11).Value = "Sin especificar"
$ws.Cells.Item(365, 12).Value = "Primera"
$ws.Cells.Item(365, 13).Value = 100
$ws.Cells.Item(365, 14).Value = 10000
$ws.Cells.Item(365, 15).Value = 10000
$ws.Cells.Item(365, 16).Value = 10000
$ws.Cells.Item(365, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(365, 18).Value = "Perú"
$ws.Cells.Item(365, 19).Value = 2500
$ws.Cells.Item(365, 20).Value = 4
